$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 2 ("Subtitle 2" placeholder): merge the two paragraphs into one,
#     drop "HTML/Bootstrap " and the "v. 1.3.0" line, and split the
#     remaining text into four runs: "Fluent Java API for " / "building " /
#     "Web " / "UI" ---
$titleShape = $s.Shapes.Item(2)
$tr = $titleShape.TextFrame.TextRange
$tr.Text = "Fluent Java API for building Web UI"

$c2 = $tr.Characters(21, 9)
$c2.Text = "building "

$c3 = $tr.Characters(30, 4)
$c3.Text = "Web "

$c4 = $tr.Characters(34, 2)
$c4.Text = "UI"

# --- Shape 3 (free-floating text box): remove the "A component of
#     Nasdanika Foundation Server" text entirely, leaving an empty
#     paragraph (endParaRPr only) ---
$fsShape = $s.Shapes.Item(3)
$fsShape.TextFrame.TextRange.Text = ""
